# Remove all old terms for MDR introduction.
# The "constants" sheet had a pair of rows describing when MDR-TB
# introduction starts/ends ("start_mdr_introduce_time" / "end_mdr_introduce_time").
# Collapse this into a single "mdr_introduce_time" row: rename the first
# row's label and delete the second (end) row entirely, letting Excel
# shift everything below it up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Row 5 held "start_mdr_introduce_time" -> rename to "mdr_introduce_time"
$ws.Range("A5").Value = "mdr_introduce_time"

# Row 6 held "end_mdr_introduce_time" (with its own description) -> delete it,
# shifting rows 7-9 up to become rows 6-8.
$ws.Rows.Item(6).Delete()

# Match the recorded selection after the edit.
[void]$ws.Range("E4").Select()
